# Word COM-interop script applying the CV content updates described by
# the commit "Language: Update fed description".

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, `
                                      $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# 1. Profile summary paragraph under "Data Scientist" heading.
Replace-Text `
    "Versatile programmer proficient in python, SQL (databases), R, etc. Skilled in Linux, OOP, data science. Teaching experience." `
    "Data Scientist experienced in data science, software engineering, and team leadership. Passionate about building innovative solutions and fostering professional growth. Adept at bridging the gap between technical and non-technical fullstack teams to achieve organizational success."

# 2. Job title line (note the non-breaking space between "Sr." and "Data").
Replace-Text `
    "Sr.$([char]160)Data Scientist, National Stress Testing, Production." `
    "Sr.$([char]160)Data Scientist, National Stress Testing Program, Production."

# 3. First bullet under Federal Reserve Bank of Minneapolis.
Replace-Text `
    "Contributed to econometric models, infrastructure, dashboards, adhoc analyses for consumption by Board of Governors (Jerome Powell)." `
    "Maintained multiple models, infrastructure, dashboards, adhoc analyses for Board of Governors."

# 4. Second bullet - liaison role.
Replace-Text `
    "Served as Production econometric modeling team$([char]8217)s technical liaison to frontend, database and sysadmin teams, ensuring smooth integration." `
    "Served as Production econometric team$([char]8217)s python lead, and liaison to database, ETL, HPC, sysadmin teams, ensuring smooth integration."

# 5. Third bullet text changes (becomes the "unit testing" bullet).
Replace-Text `
    "Technical lead in migrating Production codebase across OS, language versions, environments, while distributing and popularizing DIY automation tools to support fellow quants and economists." `
    "Proactively introduced unit testing, CI/CD, makefiles, simple-English docs to repos."

# 6. Fourth bullet text changes (becomes the "Trained economists" bullet).
Replace-Text `
    "Proactively introduced unit testing, autotesting, and autodocs to python repos." `
    "Trained economists in Python, Linux, and computer science principles, enabling self-sufficiency."

# 7. Fifth bullet text changes (becomes the "Technical lead in migrating" bullet,
#    now without the DIY-automation clause, with a reproducibility clause instead).
Replace-Text `
    "Co-prototyped a fullstack data management system using Flask and SQLite using dynamic SQL queries, HTML forms, endpoints." `
    "Technical lead in migrating Production codebase across OS, language versions, environments, enabling reproducibility."

# 8. Append a brand-new sixth bullet after the one just edited, carrying the
#    fullstack-prototype sentence that used to live in bullet 5.
$target = $null
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Technical lead in migrating Production codebase across OS, language versions, environments, enabling reproducibility.*") {
        $target = $idx
    }
}
if ($target -ne $null) {
    $p = $d.Paragraphs($target)
    $p.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs($target + 1)
    $newPara.Range.Text = "Prototyped a fullstack data management system for tracking Stress Testing operations."
} else {
    Write-Output "TARGET BULLET NOT FOUND FOR INSERTION"
}

# 9. United Health Group bullet under "outperformed competition".
Replace-Text `
    "In zero-sum ACA marketplace, United Health Group outperformed competition." `
    "In zero-sum ACA marketplace, consistently outperformed competition."

# 10. Languages line - drop "basic fluency,".
Replace-Text `
    "CEFR B2 basic fluency, but diminishing" `
    "CEFR B2, but diminishing"
